# Weekly update: insert a new data row at row 10 (pushing all subsequent
# rows down by one) with the latest week's observation for
# Feria Lagunitas de Puerto Montt - Espárragos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 10:25 down to 11:26 by inserting a new blank row at 10.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C10").Value = 'Los Lagos'
$ws.Range("D10").Value = 44498
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = 'Espárragos'
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 1600
$ws.Range("L10").Value = 1600
$ws.Range("M10").Value = 1600
$ws.Range("N10").Value = '$/kilo'
$ws.Range("O10").Value = 'Provincia de Linares'
$ws.Range("P10").Value = 1600
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'Hortaliza'
